# Compliance Officer Checklist - v1.1.1 update
#
# Adds Control 2.21 "AI Marketing Claims and Substantiation" to the
# checklist table. The new entry is inserted right after Control 2.19
# ("Customer AI Disclosure and Transparency") and before Control 3.3
# ("Compliance and Regulatory Reporting"), which (together with every
# row below it) shifts down by one row as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13 to host Control 2.21.
# Before the insert, row 13 holds "3.3 / Compliance and Regulatory Reporting".
# After the insert:
#   row 13 -> new Control 2.21 row
#   row 14 -> old row 13 (3.3 ...)
#   row 15 -> old row 14 (3.10 ...)
#   row 18 -> old row 17 (FSI Agent Governance Framework v1.1)
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the Control 2.21 details.
# Force the control id to be stored as text (like "1.7", "3.10", etc.)
# instead of being auto-converted to a number, then drop the temporary
# "@" number format so the cell keeps the default (unstyled) look used
# by its neighboring rows (2.18 / 2.19).
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "2.21"
$ws.Cells.Item(13, 1).ClearFormats()

$ws.Cells.Item(13, 2).Value = "AI Marketing Claims and Substantiation"
$ws.Cells.Item(13, 3).Value = "Not Started"

Write-Host "Inserted Control 2.21 (AI Marketing Claims and Substantiation) at row 13"
